$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 8610.5
$ws.Range("I111").Value = 8721.556
$ws.Range("J111").Value = 8277.333000000001
$ws.Range("K111").Value = 26164.668
$ws.Range("L111").Value = 24831.999
$ws.Range("M111").Value = -23097.668
$ws.Range("N111").Value = -30965.999
$ws.Range("H116").Value = 3261.625
$ws.Range("I116").Value = 2923.3333
$ws.Range("J116").Value = 3696.5715
$ws.Range("K116").Value = 2923.3333
$ws.Range("L116").Value = 3696.5715
$ws.Range("M116").Value = 518.6667000000002
$ws.Range("N116").Value = -10580.5715
$ws.Range("H129").Value = 994.3488
$ws.Range("I129").Value = 509.91666
$ws.Range("J129").Value = 1072.9054
$ws.Range("K129").Value = 1529.74998
$ws.Range("L129").Value = 3218.7162
$ws.Range("M129").Value = 3470.25002
$ws.Range("N129").Value = -13218.7162
$ws.Range("H132").Value = 2360.775
$ws.Range("I132").Value = 2304.25
$ws.Range("J132").Value = 2869.5
$ws.Range("K132").Value = 6912.75
$ws.Range("L132").Value = 8608.5
$ws.Range("M132").Value = -4382.75
$ws.Range("N132").Value = -13668.5
$ws.Range("H137").Value = 3127.0605
$ws.Range("I137").Value = 2880.3809
$ws.Range("J137").Value = 3558.75
$ws.Range("K137").Value = 8641.1427
$ws.Range("L137").Value = 10676.25
$ws.Range("M137").Value = -6091.1427
$ws.Range("N137").Value = -15776.25
$ws.Range("H138").Value = 1912.8036
$ws.Range("I138").Value = 1347.3429
$ws.Range("J138").Value = 2855.238
$ws.Range("K138").Value = 4042.0287
$ws.Range("L138").Value = 8565.714
$ws.Range("M138").Value = 1097.9713
$ws.Range("N138").Value = -18845.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3046.5652
$ws.Range("I61").Value = 2704.4375
$ws.Range("K61").Value = 2704.4375
$ws.Range("M61").Value = -2492.4375
$ws.Range("H74").Value = 1087.2059
$ws.Range("I74").Value = 749.0833
$ws.Range("J74").Value = 1898.7
$ws.Range("K74").Value = 749.0833
$ws.Range("L74").Value = 1898.7
$ws.Range("M74").Value = 124.9167
$ws.Range("N74").Value = -3646.7
$ws.Range("H77").Value = 1087.2059
$ws.Range("I77").Value = 749.0833
$ws.Range("J77").Value = 1898.7
$ws.Range("K77").Value = 3745.4165
$ws.Range("L77").Value = 9493.5
$ws.Range("M77").Value = 622.5834999999997
$ws.Range("N77").Value = -18229.5
$ws.Range("H110").Value = 1369.9166
$ws.Range("I110").Value = 1587.7778
$ws.Range("J110").Value = 716.3333
$ws.Range("K110").Value = 1587.7778
$ws.Range("L110").Value = 716.3333
$ws.Range("M110").Value = 457.2221999999999
$ws.Range("N110").Value = -4806.3333
$ws.Range("H136").Value = 3046.5652
$ws.Range("I136").Value = 2704.4375
$ws.Range("K136").Value = 8113.3125
$ws.Range("M136").Value = -5563.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1517.8572
$ws.Range("J80").Value = 233.33333
$ws.Range("L80").Value = 233.33333
$ws.Range("N80").Value = -2229.33333
$ws.Range("H83").Value = 1517.8572
$ws.Range("J83").Value = 233.33333
$ws.Range("L83").Value = 1166.66665
$ws.Range("N83").Value = -11150.66665
$ws.Range("H86").Value = 2257.8948
$ws.Range("I86").Value = 2029.1428
$ws.Range("J86").Value = 2898.4
$ws.Range("K86").Value = 2029.1428
$ws.Range("L86").Value = 2898.4
$ws.Range("M86").Value = -906.1428000000001
$ws.Range("N86").Value = -5144.4
$ws.Range("H89").Value = 2257.8948
$ws.Range("I89").Value = 2029.1428
$ws.Range("J89").Value = 2898.4
$ws.Range("K89").Value = 10145.714
$ws.Range("L89").Value = 14492
$ws.Range("M89").Value = -4529.714
$ws.Range("N89").Value = -25724
$ws.Range("H134").Value = 3057
$ws.Range("I134").Value = 2487.4285
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 7462.2855
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -4927.2855
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5865.4365
$ws.Range("I31").Value = 1277.4138
$ws.Range("J31").Value = 10982.846
$ws.Range("K31").Value = 1277.4138
$ws.Range("L31").Value = 10982.846
$ws.Range("M31").Value = -982.4138
$ws.Range("N31").Value = -11572.846
$ws.Range("H34").Value = 5865.4365
$ws.Range("I34").Value = 1277.4138
$ws.Range("J34").Value = 10982.846
$ws.Range("K34").Value = 1277.4138
$ws.Range("L34").Value = 10982.846
$ws.Range("M34").Value = -1075.4138
$ws.Range("N34").Value = -11386.846
$ws.Range("H58").Value = 1706.5217
$ws.Range("I58").Value = 1401.9166
$ws.Range("J58").Value = 2038.8182
$ws.Range("K58").Value = 1401.9166
$ws.Range("L58").Value = 2038.8182
$ws.Range("M58").Value = -1198.9166
$ws.Range("N58").Value = -2444.8182
$ws.Range("H107").Value = 805.25
$ws.Range("I107").Value = 341.66666
$ws.Range("K107").Value = 341.66666
$ws.Range("M107").Value = 1578.33334
$ws.Range("H132").Value = 1495.3889
$ws.Range("I132").Value = 1138.1428
$ws.Range("J132").Value = 1995.5333
$ws.Range("K132").Value = 3414.4284
$ws.Range("L132").Value = 5986.5999
$ws.Range("M132").Value = -884.4284000000002
$ws.Range("N132").Value = -11046.5999
$ws.Range("H136").Value = 1706.5217
$ws.Range("I136").Value = 1401.9166
$ws.Range("J136").Value = 2038.8182
$ws.Range("K136").Value = 4205.7498
$ws.Range("L136").Value = 6116.4546
$ws.Range("M136").Value = -1655.7498
$ws.Range("N136").Value = -11216.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2571.4285
$ws.Range("I19").Value = 2000
$ws.Range("J19").Value = 3333.3333
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 9999.999899999999
$ws.Range("M19").Value = -5826
$ws.Range("N19").Value = -10347.9999
$ws.Range("H38").Value = 808.95
$ws.Range("I38").Value = 75
$ws.Range("J38").Value = 1298.25
$ws.Range("K38").Value = 225
$ws.Range("L38").Value = 3894.75
$ws.Range("M38").Value = 122
$ws.Range("N38").Value = -4588.75
$ws.Range("H113").Value = 776.41174
$ws.Range("I113").Value = 557.5185
$ws.Range("J113").Value = 1022.6667
$ws.Range("K113").Value = 1672.5555
$ws.Range("L113").Value = 3068.0001
$ws.Range("M113").Value = 497.4445000000001
$ws.Range("N113").Value = -7408.0001
$ws.Range("H123").Value = 5566.25
$ws.Range("J123").Value = 5928.5713
$ws.Range("L123").Value = 17785.7139
$ws.Range("N123").Value = -22685.7139
$ws.Range("H131").Value = 952.25806
$ws.Range("I131").Value = 553.25
$ws.Range("J131").Value = 1204.2632
$ws.Range("K131").Value = 1659.75
$ws.Range("L131").Value = 3612.7896
$ws.Range("M131").Value = 3380.25
$ws.Range("N131").Value = -13692.7896
$ws.Range("H140").Value = 1899.2307
$ws.Range("J140").Value = 1893.3334
$ws.Range("L140").Value = 5680.0002
$ws.Range("N140").Value = -16040.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2679.9092
$ws.Range("I132").Value = 1693.6154
$ws.Range("J132").Value = 4104.5557
$ws.Range("K132").Value = 5080.8462
$ws.Range("L132").Value = 12313.6671
$ws.Range("M132").Value = -2550.8462
$ws.Range("N132").Value = -17373.6671
$ws.Range("H136").Value = 2018.3636
$ws.Range("I136").Value = 2901.75
$ws.Range("J136").Value = 1513.5714
$ws.Range("K136").Value = 8705.25
$ws.Range("L136").Value = 4540.7142
$ws.Range("M136").Value = -6155.25
$ws.Range("N136").Value = -9640.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2480.7368
$ws.Range("I132").Value = 2198.4167
$ws.Range("J132").Value = 2964.7144
$ws.Range("K132").Value = 6595.250100000001
$ws.Range("L132").Value = 8894.143199999999
$ws.Range("M132").Value = -4065.250100000001
$ws.Range("N132").Value = -13954.1432
$ws.Range("H136").Value = 2216.1191
$ws.Range("I136").Value = 1692.1111
$ws.Range("J136").Value = 3159.3333
$ws.Range("K136").Value = 5076.3333
$ws.Range("L136").Value = 9477.999899999999
$ws.Range("M136").Value = -2526.3333
$ws.Range("N136").Value = -14577.9999
